# Generate Report for Handoff
# - Set Priority ("ht") for the rows whose handoff xliff was (re)generated,
#   on both the zh-cn and de-de sheets.
# - Refresh the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#   timestamps for those same rows on Overview, zh-cn and de-de.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 11, 12, 13, 14)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

foreach ($r in $rows) {
    # Priority column (E) on the per-language sheets: blank -> "ht"
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 5).Value = "ht"

    # Latest Handoff Datetime (column H) on the per-language sheets
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-30 12:25:27"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-30 12:25:32"

    # Latest HO Xliff Generate Date (column G) on the Overview sheet
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-30 12:25:32"
}
